$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '51.605.23'
$ws.Cells.Item(2, 5).Value = '  +1.09%  '
# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.029.35'
$ws.Cells.Item(3, 5).Value = '  +2.46%  '
# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '382.12'
$ws.Cells.Item(5, 5).Value = '  +0.41%  '
# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '102.61'
$ws.Cells.Item(6, 5).Value = '  +0.49%  '
# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.03%  '
# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.590'
$ws.Cells.Item(9, 5).Value = '  +0.35%  '
# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '36.69'
$ws.Cells.Item(10, 5).Value = '  +0.51%  '
# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.138'
$ws.Cells.Item(11, 5).Value = '  +0.02%  '
# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0861'
$ws.Cells.Item(12, 5).Value = '  +1.11%  '
# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '3.499.50'
$ws.Cells.Item(13, 5).Value = '  +2.33%  '
# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '18.51'
$ws.Cells.Item(14, 5).Value = '  +0.75%  '
# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.73'
$ws.Cells.Item(15, 5).Value = '  -0.08%  '
# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.026.35'
$ws.Cells.Item(16, 5).Value = '  +2.10%  '
# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.973'
$ws.Cells.Item(17, 5).Value = '  -3.17%  '
# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '10.66'
$ws.Cells.Item(18, 5).Value = '  -13.77%  '
# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '51.600.22'
$ws.Cells.Item(19, 5).Value = '  +0.96%  '
# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '3.08'
$ws.Cells.Item(20, 5).Value = '  -0.24%  '
# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '12.47'
$ws.Cells.Item(21, 5).Value = '  +0.84%  '
# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.0₃0962'
$ws.Cells.Item(22, 5).Value = '  +0.31%  '
# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '69.95'
$ws.Cells.Item(23, 5).Value = '  +0.32%  '
# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '267.33'
$ws.Cells.Item(24, 5).Value = '  -0.66%  '
# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '3.16'
$ws.Cells.Item(25, 5).Value = '  -5.26%  '
# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '8.44'
$ws.Cells.Item(26, 5).Value = '  +6.51%  '
# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '7.52'
$ws.Cells.Item(27, 5).Value = '  +6.73%  '
# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.172'
$ws.Cells.Item(28, 5).Value = '  +3.49%  '
# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.07%  '
# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '26.20'
$ws.Cells.Item(30, 5).Value = '  +1.18%  '
# Row 31
$ws.Cells.Item(31, 5).Value = '  -0.97%  '
# Row 32
$ws.Cells.Item(32, 5).Value = '  -1.58%  '
# Row 33
$ws.Cells.Item(33, 2).Value = 'Toncoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.08'
$ws.Cells.Item(33, 5).Value = '  -2.30%  '
# Row 34
$ws.Cells.Item(34, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '34.12'
$ws.Cells.Item(34, 5).Value = '  -0.70%  '
# Row 35
$ws.Cells.Item(35, 2).Value = 'OKB'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '50.54'
$ws.Cells.Item(35, 5).Value = '  -1.19%  '
# Row 36
$ws.Cells.Item(36, 2).Value = 'VeChain'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0448'
$ws.Cells.Item(36, 5).Value = '  +2.75%  '
# Row 37
$ws.Cells.Item(37, 5).Value = '  -0.12%  '
# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.30'
$ws.Cells.Item(38, 5).Value = '  +1.76%  '
# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.298'
$ws.Cells.Item(39, 5).Value = '  +9.09%  '
# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '17.05'
$ws.Cells.Item(40, 5).Value = '  +2.01%  '
# Row 41
$ws.Cells.Item(41, 2).Value = 'ARBITRUM'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.86'
$ws.Cells.Item(41, 5).Value = '  +1.37%  '
# Row 42
$ws.Cells.Item(42, 2).Value = 'Monero'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '127.48'
$ws.Cells.Item(42, 5).Value = '  +2.31%  '
# Row 43
$ws.Cells.Item(43, 5).Value = '  -0.77%  '
# Row 44
$ws.Cells.Item(44, 5).Value = '  +0.74%  '
# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.71'
$ws.Cells.Item(45, 5).Value = '  +4.41%  '
# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '21.56'
$ws.Cells.Item(46, 5).Value = '  -0.30%  '
# Row 47
$ws.Cells.Item(47, 5).Value = '  +2.66%  '
# Row 48
$ws.Cells.Item(48, 5).Value = '  +2.52%  '
# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.022.96'
$ws.Cells.Item(49, 5).Value = '  -2.12%  '
# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '3.326.48'
$ws.Cells.Item(50, 5).Value = '  +2.39%  '
# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.514'
$ws.Cells.Item(51, 5).Value = '  +5.30%  '

Write-Output "Update complete"